$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, matching the target diff
$updates = @{
    2  = @{ 'E' = 3; 'F' = 1; 'G' = 0.1346003333333333; 'H' = 0.403801; 'I' = 0.009651054304565105; 'J' = 0.009651054304565105;
            'M' = 145.520762; 'N' = 436.562286; 'O' = 0.962605808965774; 'P' = 0.9626058089657743;
            'Q' = 19.58714307212067; 'R' = 176.284287649086; 'S' = 0.009290160936218508; 'T' = 0.009290160936218511 }
    3  = @{ 'E' = 3; 'F' = 1; 'G' = 0.1346003333333333; 'H' = 0.403801; 'I' = 0.009651054304565105; 'J' = 0.009651054304565105;
            'O' = 0.02468611455667964; 'P' = 0.02468611455667964;
            'Q' = 0.5023140866311112; 'R' = 4.52082677968; 'S' = 0.0002382470321552303; 'T' = 0.0002382470321552303 }
    4  = @{ 'E' = 3; 'F' = 1; 'G' = 0.1346003333333333; 'H' = 0.403801; 'I' = 0.009651054304565105; 'J' = 0.009651054304565105;
            'M' = 1.921128; 'N' = 5.763384; 'O' = 0.01270807647754621; 'P' = 0.01270807647754621;
            'Q' = 0.258584469176; 'R' = 2.327260222584; 'S' = 0.0001226463361913649; 'T' = 0.0001226463361913649 }
    5  = @{ 'I' = 0.8124788779145131; 'J' = 0.8124788779145132;
            'M' = 145.520762; 'N' = 436.562286; 'O' = 0.962605808965774; 'P' = 0.9626058089657743;
            'Q' = 1648.953525964514; 'R' = 14840.58173368062; 'S' = 0.7820968875425042; 'T' = 0.7820968875425046 }
    6  = @{ 'I' = 0.8124788779145131; 'J' = 0.8124788779145132;
            'O' = 0.02468611455667964; 'P' = 0.02468611455667964;
            'S' = 0.0200569466550802; 'T' = 0.0200569466550802 }
    7  = @{ 'I' = 0.8124788779145131; 'J' = 0.8124788779145132;
            'M' = 1.921128; 'N' = 5.763384; 'O' = 0.01270807647754621; 'P' = 0.01270807647754621;
            'Q' = 21.769064055816; 'R' = 195.921576502344; 'S' = 0.01032504371692856; 'T' = 0.01032504371692856 }
    8  = @{ 'G' = 2.4807; 'H' = 7.4421; 'I' = 0.1778700677809217; 'J' = 0.1778700677809217;
            'M' = 145.520762; 'N' = 436.562286; 'O' = 0.962605808965774; 'P' = 0.9626058089657743;
            'Q' = 360.9933542934; 'R' = 3248.9401886406; 'S' = 0.1712187604870512; 'T' = 0.1712187604870513 }
    9  = @{ 'G' = 2.4807; 'H' = 7.4421; 'I' = 0.1778700677809217; 'J' = 0.1778700677809217;
            'O' = 0.02468611455667964; 'P' = 0.02468611455667964;
            'Q' = 9.257707792; 'R' = 83.31937012799999; 'S' = 0.004390920869444205; 'T' = 0.004390920869444205 }
    10 = @{ 'G' = 2.4807; 'H' = 7.4421; 'I' = 0.1778700677809217; 'J' = 0.1778700677809217;
            'M' = 1.921128; 'N' = 5.763384; 'O' = 0.01270807647754621; 'P' = 0.01270807647754621;
            'Q' = 4.765742229600001; 'R' = 42.8916800664; 'S' = 0.002260386424426281; 'T' = 0.002260386424426281 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

$wb.Save()
